# Germany Landesliga workbook update (02-05-2024 20:28)
#
# The underlying source data was re-sorted, which caused three pairs of
# match rows (which share the same kick-off date) to swap places with
# each other. Swap the full records (every column except the running
# index in column A, which stays positional) for each pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($sheet, [int]$row1, [int]$row2) {
    $range1 = $sheet.Range("B$row1`:AB$row1")
    $range2 = $sheet.Range("B$row2`:AB$row2")
    $values1 = $range1.Value2
    $values2 = $range2.Value2
    $range1.Value = $values2
    $range2.Value = $values1
}

Swap-Rows $ws 4 5
Swap-Rows $ws 12 13
Swap-Rows $ws 85 86
